$wb = $excel.ActiveWorkbook

$items = $wb.Worksheets.Item("items")
$models = $wb.Worksheets.Item("models")
$pricing = $wb.Worksheets.Item("book pricing")

# --- sheet "items": rename Banner -> Banners -----------------------------
$items.Range("A2").Value = "Banners"

# --- sheet "models": fix pricing / remove the stray Book row -------------

# The "Book / Hardcover / Hardcover Book / 300" row is bogus -> delete it.
$models.Rows.Item(4).Delete() | Out-Null

# Row 2 ("Banner" header row) -> rename + correct price
$models.Range("A2").Value = "Banners"
$models.Range("D2").Value = 50

# Correct the (now shifted) "Books" binding rows' prices
$models.Range("D5").Value = 1.2
$models.Range("D6").Value = 1
$models.Range("D7").Value = 0.8
$models.Range("D8").Value = 1.2

# New rows describing the "Banners" model variants
$models.Range("A13").Value = "Banners"
$models.Range("B13").Value = "Sticker"
$models.Range("D13").Value = 60
$models.Range("D13").VerticalAlignment = -4108

$models.Range("A14").Value = "Banners"
$models.Range("B14").Value = "Poster"
$models.Range("D14").Value = 80
$models.Range("D14").VerticalAlignment = -4108

$models.Range("A15").Value = "Banners"
$models.Range("B15").Value = "Reflective"
$models.Range("D15").Value = 180
$models.Range("D15").VerticalAlignment = -4108

# --- sheet "book pricing": only a cosmetic column width tweak ------------
$pricing.Columns.Item(1).ColumnWidth = 13.43

# --- selection / active-sheet bookkeeping ---------------------------------
$items.Range("A3").Select() | Out-Null
$models.Activate() | Out-Null
$models.Range("D8").Select() | Out-Null
